$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.437.94"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "3.071.38"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'522.42"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'140.14"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.069.12"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("D13").Value = "3.599.34"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "'25.30"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "57.486.49"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "3.069.79"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'6.05"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").Value = "'12.61"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "'7.90"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("D22").Value = "'338.49"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'0.507"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "'66.79"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "0.0₃0905"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'6.33"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").Value = "'7.18"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("D33").Value = "'20.81"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("D35").Value = "'157.76"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Value = "'25.77"
$ws.Range("E38").Value = "  -5.82%  "
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").Value = "'0.0661"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("E41").Value = "  +12.62%  "
$ws.Range("D42").Value = "'3.96"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "'0.679"
$ws.Range("E43").Value = "  +3.53%  "
$ws.Range("D44").Value = "3.108.40"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "'36.70"
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "2.269.10"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E49").Value = "  +4.44%  "
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").Value = "'20.39"
$ws.Range("E51").Value = "  -0.64%  "